$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.946.86"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").Value = "2.574.91"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.55"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.31"
$ws.Range("E6").Value = "  +3.41%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +0.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.13"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.56"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").Value = "2.970.87"
$ws.Range("E13").Value = "  +2.03%  "

$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").Value = "2.631.92"
$ws.Range("E15").Value = "  +3.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.63"
$ws.Range("E16").Value = "  +2.21%  "

$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").Value = "43.013.17"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("E19").Value = "  +2.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("E20").Value = "  -2.09%  "

$ws.Range("D21").Value = "0.0₃0970"
$ws.Range("E21").Value = "  +0.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.38"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.12"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.06"
$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  -1.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.58"
$ws.Range("E29").Value = "  -2.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.29"
$ws.Range("E30").Value = "  -0.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.85"
$ws.Range("E31").Value = "  -2.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.54"
$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.44"
$ws.Range("E33").Value = "  +4.03%  "

$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0804"
$ws.Range("E35").Value = "  +2.94%  "

$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.79"
$ws.Range("E37").Value = "  -2.66%  "

$ws.Range("B38").Value = "ApeXProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.52"
$ws.Range("E38").Value = "  +9.55%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.112"
$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.11"
$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("E42").Value = "  +5.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0305"
$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("E45").Value = "  -2.56%  "

$ws.Range("D46").Value = "2.007.16"
$ws.Range("E46").Value = "  -1.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.90"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("D48").Value = "2.822.72"
$ws.Range("E48").Value = "  +2.13%  "

$ws.Range("E49").Value = "  +2.04%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.06"
$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.03"
$ws.Range("E51").Value = "  -3.03%  "

